# The document has two logos embedded as inline pictures in the page
# headers (BTEC logo) and footers (Pearson logo). Each logo's drawing
# object carries a "name" (exposed on the Word object model as
# InlineShape.Name, which maps to the OOXML <wp:docPr name="..."/>
# attribute). This script swaps those names:
#   image1.png -> image2.png   (Pearson logo, in the footers)
#   image2.jpg -> image1.jpg   (BTEC logo, in the headers)

$d = $word.ActiveDocument

function Update-InlineShapeName($shape) {
    $xml = $shape.WordOpenXML
    if ($xml -match '<wp:docPr[^>]*\bname="([^"]*)"') {
        $current = $matches[1]
        $newName = $null
        if ($current -eq "image1.png") {
            $newName = "image2.png"
        } elseif ($current -eq "image2.jpg") {
            $newName = "image1.jpg"
        }
        if ($newName -ne $null -and $newName -ne $current) {
            $shape.Name = $newName
        }
    }
}

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $h = $sec.Headers.Item($i)
        if ($h.Exists) {
            $shapes = $h.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                Update-InlineShapeName $shapes.Item($j)
            }
        }
    }
    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $f = $sec.Footers.Item($i)
        if ($f.Exists) {
            $shapes = $f.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                Update-InlineShapeName $shapes.Item($j)
            }
        }
    }
}

Write-Output "Renamed inline shapes in headers/footers"
